$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row after row 10 (old SEPARATOR row) to make room for the
# "rating" row that will replace the old row 10, pushing the SEPARATOR
# (and everything below it) down by one row.
$ws.Rows.Item(11).Insert()

# The old SEPARATOR that used to live in row 10 now belongs in row 11.
$ws.Range("A11").Value = "SEPARATOR"

# Row 10 becomes the new "rating" entry (was "SEPARATOR").
$ws.Range("A10").Value = "rating"
$ws.Range("B10").Value = 6.66
$ws.Range("B10").NumberFormat = "0.00"

# Append a new "rating" row at the very end of the sheet (row 20).
$ws.Range("A20").Value = "rating"
$ws.Range("B20").Value = 6.66

# Leave the selection where the author last left it.
$null = $ws.Range("D18").Select()
